$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33, shifting existing rows 33-115 down to 34-116
$ws.Rows("33:33").Insert()

# Populate the newly inserted row 33 with the new weekly record
$ws.Range("A33").Value = 7
$ws.Range("B33").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C33").Value = 'Ñuble'
$ws.Range("D33").Value = 44910
$ws.Range("D33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E33").Value = 16
$ws.Range("F33").Value = 100112021
$ws.Range("G33").Value = 'Ají'
$ws.Range("H33").Value = 'Americana (o)'
$ws.Range("I33").Value = 'Primera'
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 15000
$ws.Range("L33").Value = 15000
$ws.Range("M33").Value = 15000
$ws.Range("N33").Value = '$/caja 15 kilos'
$ws.Range("O33").Value = 'Región del Maule'
$ws.Range("P33").Value = 1000
$ws.Range("Q33").Value = 15
$ws.Range("R33").Value = 'Hortaliza'
